$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '61.836.11'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '3.407.58'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").Value = '409.21'
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").Value = '129.04'
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("D7").Value = '0.620'
$ws.Range("E7").Value = '  -1.94%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.726'
$ws.Range("E9").Value = '  -2.47%  '
$ws.Range("D10").Value = '0.134'
$ws.Range("E10").Value = '  -6.92%  '
$ws.Range("D11").Value = '42.66'
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '9.13'
$ws.Range("E12").Value = '  +1.76%  '
$ws.Range("D13").Value = '3.944.21'
$ws.Range("E13").Value = '  -0.73%  '
$ws.Range("D14").Value = '0.140'
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").Value = '0.0000209'
$ws.Range("E15").Value = '  -3.97%  '
$ws.Range("D16").Value = '20.38'
$ws.Range("E16").Value = '  -2.99%  '
$ws.Range("D17").Value = '3.419.26'
$ws.Range("E17").Value = '  +1.80%  '
$ws.Range("D18").Value = '1.08'
$ws.Range("E18").Value = '  +1.10%  '
$ws.Range("D19").Value = '12.25'
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = '61.840.58'
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("D21").Value = '484.96'
$ws.Range("E21").Value = '  +20.93%  '
$ws.Range("D22").Value = '90.06'
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("D23").Value = '3.26'
$ws.Range("E23").Value = '  +2.32%  '
$ws.Range("D24").Value = '13.18'
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("D25").Value = '3.28'
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("D26").Value = '9.51'
$ws.Range("E26").Value = '  +10.35%  '
$ws.Range("D27").Value = '33.22'
$ws.Range("E27").Value = '  +1.83%  '
$ws.Range("D28").Value = '4.81'
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '7.93'
$ws.Range("E29").Value = '  +4.36%  '
$ws.Range("D30").Value = '2.65'
$ws.Range("E30").Value = '  -3.07%  '
$ws.Range("D31").Value = '11.84'
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("D32").Value = '0.168'
$ws.Range("E32").Value = '  -2.10%  '
$ws.Range("D33").Value = '0.113'
$ws.Range("E33").Value = '  -6.27%  '
$ws.Range("D34").Value = '41.09'
$ws.Range("E34").Value = '  -5.03%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.70%  '
$ws.Range("D36").Value = '56.94'
$ws.Range("E36").Value = '  +5.41%  '
$ws.Range("D37").Value = '0.0487'
$ws.Range("E37").Value = '  -2.16%  '
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("D39").Value = '3.00'
$ws.Range("E39").Value = '  +2.99%  '
$ws.Range("D40").Value = '0.325'
$ws.Range("E40").Value = '  +4.48%  '
$ws.Range("D41").Value = '147.81'
$ws.Range("E41").Value = '  +4.66%  '
$ws.Range("D42").Value = '0.134'
$ws.Range("E42").Value = '  +0.57%  '
$ws.Range("D43").Value = '3.33'
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("D44").Value = '2.08'
$ws.Range("E44").Value = '  +5.16%  '
$ws.Range("D45").Value = '2.59'
$ws.Range("E45").Value = '  +7.47%  '
$ws.Range("D46").Value = '4.22'
$ws.Range("E46").Value = '  +3.53%  '
$ws.Range("D47").Value = '2.36'
$ws.Range("E47").Value = '  +18.85%  '
$ws.Range("D48").Value = '16.56'
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("D49").Value = '22.04'
$ws.Range("E49").Value = '  +1.02%  '
$ws.Range("D50").Value = '112.86'
$ws.Range("E50").Value = '  +14.20%  '
$ws.Range("D51").Value = '0.140'
$ws.Range("E51").Value = '  +5.48%  '
